$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.032805416042343
$ws.Range("D2").Value = 1.034554572569423
$ws.Range("E2").Value = 1.042269083816134
$ws.Range("F2").Value = 1.053454673562189
$ws.Range("I2").Value = 1.033340783166479
$ws.Range("J2").Value = 1.037933617019257
$ws.Range("K2").Value = 1.037353712205781
$ws.Range("L2").Value = 1.045046248859688
$ws.Range("M2").Value = 1.056200596468772
$ws.Range("N2").Value = 1.016555306158375

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.033927637738422
$ws.Range("D3").Value = 1.035352976578142
$ws.Range("E3").Value = 1.043287129538982
$ws.Range("F3").Value = 1.054607117566863
$ws.Range("I3").Value = 1.033537193330002
$ws.Range("J3").Value = 1.038697410743754
$ws.Range("K3").Value = 1.037961356930927
$ws.Range("L3").Value = 1.045874548556565
$ws.Range("M3").Value = 1.057165219297967
$ws.Range("N3").Value = 1.016813969316331

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.034653838272495
$ws.Range("D4").Value = 1.035869423301159
$ws.Range("E4").Value = 1.043946271417894
$ws.Range("F4").Value = 1.055353362840548
$ws.Range("I4").Value = 1.033662883921198
$ws.Range("J4").Value = 1.039191153557035
$ws.Range("K4").Value = 1.038353718259019
$ws.Range("L4").Value = 1.046410305520171
$ws.Range("M4").Value = 1.057789350348989
$ws.Range("N4").Value = 1.016981037712857

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.034959145246177
$ws.Range("D5").Value = 1.03608649549964
$ws.Range("E5").Value = 1.044223470093225
$ws.Range("F5").Value = 1.055667213195866
$ws.Range("I5").Value = 1.03371538897224
$ws.Range("J5").Value = 1.039398607591851
$ws.Range("K5").Value = 1.038518469246502
$ws.Range("L5").Value = 1.046635488003116
$ws.Range("M5").Value = 1.058051724613011
$ws.Range("N5").Value = 1.017051200447921

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.035010408372245
$ws.Range("D6").Value = 1.036122940404062
$ws.Range("E6").Value = 1.044270018569346
$ws.Range("F6").Value = 1.05571991761561
$ws.Range("I6").Value = 1.033724185136966
$ws.Range("J6").Value = 1.039433433271521
$ws.Range("K6").Value = 1.038546120065156
$ws.Range("L6").Value = 1.046673294219591
$ws.Range("M6").Value = 1.058095777800828
$ws.Range("N6").Value = 1.017062976808292

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.034657917750047
$ws.Range("D7").Value = 1.035872323995972
$ws.Range("E7").Value = 1.043949974985069
$ws.Range("F7").Value = 1.055357556017719
$ws.Range("I7").Value = 1.033663586813819
$ws.Range("J7").Value = 1.039193926021445
$ws.Range("K7").Value = 1.038355920446795
$ws.Range("L7").Value = 1.046413314614831
$ws.Range("M7").Value = 1.057792856247472
$ws.Range("N7").Value = 1.016981975516948

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.033184666628012
$ws.Range("D8").Value = 1.034824432782416
$ws.Range("E8").Value = 1.042613054737378
$ws.Range("F8").Value = 1.053844036411037
$ws.Range("I8").Value = 1.033407450684949
$ws.Range("J8").Value = 1.038191844735634
$ws.Range("K8").Value = 1.037559239437599
$ws.Range("L8").Value = 1.045326219833951
$ws.Range("M8").Value = 1.056526604591445
$ws.Range("N8").Value = 1.016642785491787

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.03058894343929
$ws.Range("D9").Value = 1.032976595726709
$ws.Range("E9").Value = 1.040260275584839
$ws.Range("F9").Value = 1.051181120023462
$ws.Range("I9").Value = 1.03294538995349
$ws.Range("J9").Value = 1.036422345106051
$ws.Range("K9").Value = 1.036149067884325
$ws.Range("L9").Value = 1.043409018141187
$ws.Range("M9").Value = 1.054294957528947
$ws.Range("N9").Value = 1.016042763317279

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.028858634624514
$ws.Range("D10").Value = 1.031743834380871
$ws.Range("E10").Value = 1.038693793872965
$ws.Range("F10").Value = 1.049408575708672
$ws.Range("I10").Value = 1.032630150963067
$ws.Range("J10").Value = 1.035240176185957
$ws.Range("K10").Value = 1.035204708976158
$ws.Range("L10").Value = 1.042129795368322
$ws.Range("M10").Value = 1.052806941273512
$ws.Range("N10").Value = 1.015641186362753

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.028109420739306
$ws.Range("D11").Value = 1.031209831502662
$ws.Range("E11").Value = 1.038015970230963
$ws.Range("F11").Value = 1.04864168681159
$ws.Range("I11").Value = 1.032491942605356
$ws.Range("J11").Value = 1.034727686988734
$ws.Range("K11").Value = 1.034794784324038
$ws.Range("L11").Value = 1.041575616303104
$ws.Range("M11").Value = 1.052162549614698
$ws.Range("N11").Value = 1.015466928626736

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.027831131383234
$ws.Range("D12").Value = 1.031011447771294
$ws.Range("E12").Value = 1.037764267084053
$ws.Range("F12").Value = 1.048356924833739
$ws.Range("I12").Value = 1.032440349312257
$ws.Range("J12").Value = 1.034537234806012
$ws.Range("K12").Value = 1.034642368059859
$ws.Range("L12").Value = 1.041369728998852
$ws.Range("M12").Value = 1.051923182727753
$ws.Range("N12").Value = 1.015402145579478

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.027890825320992
$ws.Range("D13").Value = 1.031054003168865
$ws.Range("E13").Value = 1.037818255068598
$ws.Range("F13").Value = 1.048418002969029
$ws.Range("I13").Value = 1.032451427858087
$ws.Range("J13").Value = 1.034578091564397
$ws.Range("K13").Value = 1.034675068745749
$ws.Range("L13").Value = 1.04141389436455
$ws.Range("M13").Value = 1.051974528249329
$ws.Range("N13").Value = 1.015416044296172

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.028086417219376
$ws.Range("D14").Value = 1.031193433679047
$ws.Range("E14").Value = 1.037995162929518
$ws.Range("F14").Value = 1.048618146369318
$ws.Range("I14").Value = 1.03248768312118
$ws.Range("J14").Value = 1.034711946008031
$ws.Range("K14").Value = 1.034782188649379
$ws.Range("L14").Value = 1.041558598423234
$ws.Range("M14").Value = 1.052142763703139
$ws.Range("N14").Value = 1.015461574781731

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.028206928079128
$ws.Range("D15").Value = 1.031279337276187
$ws.Range("E15").Value = 1.038104171063148
$ws.Range("F15").Value = 1.048741473854485
$ws.Range("I15").Value = 1.032509987189924
$ws.Range("J15").Value = 1.034794406094708
$ws.Range("K15").Value = 1.034848168613288
$ws.Range("L15").Value = 1.041647749988157
$ws.Range("M15").Value = 1.052246417628587
$ws.Range("N15").Value = 1.015489620199523

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.028908357831874
$ws.Range("D16").Value = 1.031779269998947
$ws.Range("E16").Value = 1.038738788735502
$ws.Range("F16").Value = 1.049459484914057
$ws.Range("I16").Value = 1.032639287417153
$ws.Range("J16").Value = 1.035274175683685
$ws.Range("K16").Value = 1.035231893015605
$ws.Range("L16").Value = 1.04216656874595
$ws.Range("M16").Value = 1.052849705891599
$ws.Range("N16").Value = 1.015652743428486

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.029348351059594
$ws.Range("D17").Value = 1.03209280870126
$ws.Range("E17").Value = 1.039136994083454
$ws.Range("F17").Value = 1.04991004371489
$ws.Range("I17").Value = 1.032719936823752
$ws.Range("J17").Value = 1.035574960953078
$ws.Range("K17").Value = 1.035472322446212
$ws.Range("L17").Value = 1.042491938246249
$ws.Range("M17").Value = 1.053228113374853
$ws.Range("N17").Value = 1.015754966628978

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.029604993792824
$ws.Range("D18").Value = 1.032275670340162
$ws.Range("E18").Value = 1.039369306215056
$ws.Range("F18").Value = 1.050172908301099
$ws.Range("I18").Value = 1.032766813479279
$ws.Range("J18").Value = 1.035750345829592
$ws.Range("K18").Value = 1.035612463324249
$ws.Range("L18").Value = 1.042681695019562
$ws.Range("M18").Value = 1.053448825259733
$ws.Range("N18").Value = 1.015814555776721

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.029692502728573
$ws.Range("D19").Value = 1.032338018006442
$ws.Range("E19").Value = 1.039448526409392
$ws.Range("F19").Value = 1.050262548741687
$ws.Range("I19").Value = 1.032782769257405
$ws.Range("J19").Value = 1.035810137700364
$ws.Range("K19").Value = 1.035660231200551
$ws.Range("L19").Value = 1.042746392819368
$ws.Range("M19").Value = 1.053524081170816
$ws.Range("N19").Value = 1.015834868045368

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.029301143738377
$ws.Range("D20").Value = 1.032059171052249
$ws.Range("E20").Value = 1.039094265694877
$ws.Range("F20").Value = 1.049861696687456
$ws.Range("I20").Value = 1.032711300949331
$ws.Range("J20").Value = 1.035542695550714
$ws.Range("K20").Value = 1.035446536733483
$ws.Range("L20").Value = 1.042457031844379
$ws.Range("M20").Value = 1.053187514539691
$ws.Range("N20").Value = 1.015744002762139

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.02802882023325
$ws.Range("D21").Value = 1.031152375731659
$ws.Range("E21").Value = 1.037943066024169
$ws.Range("F21").Value = 1.048559206517055
$ws.Range("I21").Value = 1.032477013931846
$ws.Range("J21").Value = 1.034672531722416
$ws.Range("K21").Value = 1.034750648704336
$ws.Range("L21").Value = 1.041515987806387
$ws.Range("M21").Value = 1.052093222866469
$ws.Range("N21").Value = 1.015448168733766

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.027228871178122
$ws.Range("D22").Value = 1.030582056572546
$ws.Range("E22").Value = 1.037219670160321
$ws.Range("F22").Value = 1.047740826873997
$ws.Range("I22").Value = 1.03232822404597
$ws.Range("J22").Value = 1.034124899208959
$ws.Range("K22").Value = 1.034312236823816
$ws.Range("L22").Value = 1.040924081415303
$ws.Range("M22").Value = 1.051405133599678
$ws.Range("N22").Value = 1.015261842321871

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.027652938628012
$ws.Range("D23").Value = 1.030884410660365
$ws.Range("E23").Value = 1.037603117344726
$ws.Range("F23").Value = 1.048174613705914
$ws.Range("I23").Value = 1.032407241070673
$ws.Range("J23").Value = 1.034415259554495
$ws.Range("K23").Value = 1.034544730653955
$ws.Range("L23").Value = 1.041237884564223
$ws.Range("M23").Value = 1.051769909033847
$ws.Range("N23").Value = 1.015360648184825

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.029322474687878
$ws.Range("D24").Value = 1.032074370520172
$ws.Range("E24").Value = 1.039113572673073
$ws.Range("F24").Value = 1.049883542438527
$ws.Range("I24").Value = 1.032715203638337
$ws.Range("J24").Value = 1.03555727507646
$ws.Range("K24").Value = 1.035458188488383
$ws.Range("L24").Value = 1.042472804624223
$ws.Range("M24").Value = 1.053205859426512
$ws.Range("N24").Value = 1.015748956972551

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.031259965827035
$ws.Range("D25").Value = 1.033454460450133
$ws.Range("E25").Value = 1.040868165516526
$ws.Range("F25").Value = 1.051869064342855
$ws.Range("I25").Value = 1.033066112983435
$ws.Range("J25").Value = 1.036880243382014
$ws.Range("K25").Value = 1.036514379389026
$ws.Range("L25").Value = 1.043904851720169
$ws.Range("M25").Value = 1.054871934926294
$ws.Range("N25").Value = 1.016198158896076
